# Re-sort the worksheet tabs: put "总计" (summary) before "2022-Q2" (detail).
# Before: [2022-Q2, 总计]   ->   After: [总计, 2022-Q2]
# No cell data changes - this is purely a sheet reorder.
$wb = $excel.ActiveWorkbook

# Move "总计" so it lands immediately before "2022-Q2" -> becomes the first tab.
$wb.Worksheets.Item("总计").Move($wb.Worksheets.Item("2022-Q2"))

# Keep "2022-Q2" as the active/selected sheet, same as before the reorder.
# (re-fetch by name rather than reusing a pre-move object reference)
$wb.Worksheets.Item("2022-Q2").Activate()
